# Update "想去人数" (number of people wanting to attend) counts for several
# events on both the "展览" sheet and the "全部类型" sheet, which mirror the
# same underlying data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 2347
    $ws.Range("F4").Value = 422
    $ws.Range("F5").Value = 83
    $ws.Range("F6").Value = 6479
    $ws.Range("F7").Value = 332
}
